$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 16.66879653930664
$ws.Range("C3").Value = 16.03007316589355
$ws.Range("C4").Value = 16.1888599395752
$ws.Range("C5").Value = 16.04318618774414
$ws.Range("C6").Value = 16.18576049804688
